# Auto-generated edit script applying the BOQ row 8-20 content/number
# updates and the recomputed Grand Total / Net Payable figures (rows 22 & 24)
# described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("A8").Value = ''''
$ws.Range("A8").Style = "Normal"
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = '''1.0'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = '''0.00'
$ws.Range("G8").Style = "Normal"

# --- Row 9 ---
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = '''3'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F9").Value = 472
$ws.Range("G9").Value = '''472.00'
$ws.Range("G9").Style = "Normal"

# --- Row 10 ---
$ws.Range("A10").Value = 'P. point'
$ws.Range("C10").Value = 87
$ws.Range("D10").Value = '''4'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F10").Value = 662
$ws.Range("G10").Value = '''57594.00'
$ws.Range("G10").Style = "Normal"

# --- Row 11 ---
$ws.Range("C11").Value = 29
$ws.Range("G11").Value = '''1450.00'
$ws.Range("G11").Style = "Normal"

# --- Row 12 ---
$ws.Range("C12").Value = 76
$ws.Range("D12").Value = '''8.0'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 30
$ws.Range("G12").Value = '''2280.00'
$ws.Range("G12").Style = "Normal"

# --- Row 13 ---
$ws.Range("C13").Value = 15
$ws.Range("D13").Value = '''9.0'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 219
$ws.Range("G13").Value = '''3285.00'
$ws.Range("G13").Style = "Normal"

# --- Row 14 ---
$ws.Range("A14").Value = 'Each'
$ws.Range("C14").Value = 24
$ws.Range("D14").Value = '''10.0'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F14").Value = 303
$ws.Range("G14").Value = '''7272.00'
$ws.Range("G14").Style = "Normal"

# --- Row 15 ---
$ws.Range("A15").Value = 'R. mtr.'
$ws.Range("C15").Value = 72
$ws.Range("D15").Value = '''17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '25 mm'
$ws.Range("F15").Value = 56
$ws.Range("G15").Value = '''4032.00'
$ws.Range("G15").Style = "Normal"

# --- Row 16 ---
$ws.Range("A16").Value = ''''
$ws.Range("A16").Style = "Normal"
$ws.Range("C16").Value = 76
$ws.Range("D16").Value = '''17.0'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = '''0.00'
$ws.Range("G16").Style = "Normal"

# --- Row 17 ---
$ws.Range("A17").Value = ''''
$ws.Range("A17").Style = "Normal"
$ws.Range("C17").Value = 44
$ws.Range("D17").Value = '''29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = '''0.00'
$ws.Range("G17").Style = "Normal"

# --- Row 18 ---
$ws.Range("C18").Value = 80
$ws.Range("D18").Value = '''34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'

# --- Row 19 ---
$ws.Range("A19").Value = ''''
$ws.Range("A19").Style = "Normal"
$ws.Range("C19").Value = 55
$ws.Range("D19").Value = '''36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = 'Total'
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = '''0.00'
$ws.Range("G19").Style = "Normal"

# --- Row 20 ---
$ws.Range("A20").Value = '%'
$ws.Range("C20").Value = 51
$ws.Range("D20").Value = '''37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = 'Add Tender Premium '

# --- Row 22 ---
$ws.Range("G22").Value = '''76385.00'
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value = '''76385.00'
$ws.Range("H22").Style = "Normal"

# --- Row 24 ---
$ws.Range("G24").Value = '''76385.00'
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value = '''76385.00'
$ws.Range("H24").Style = "Normal"

